# Update the statsmodels OLS summary text embedded in column B, row 2 of
# each "backward elimination step" worksheet: the run date moved from
# Saturday 28 Dec 2019 to Sunday 29 Dec 2019, and the wall-clock time of
# the run was refreshed to match the new save.
$wb = $excel.ActiveWorkbook

# New time-of-day per worksheet (tab order 1..29), taken from the target
# diff: the first 23 sheets land on 16:11:27, the remaining 6 on 16:11:28.
$newTimes = @(
    "16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27",
    "16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27",
    "16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27","16:11:27",
    "16:11:28","16:11:28","16:11:28","16:11:28","16:11:28","16:11:28"
)

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value2.ToString()

    $newTime = $newTimes[$i - 1]

    $updated = $text.Replace("Date:                Sat, 28 Dec 2019", "Date:                Sun, 29 Dec 2019")
    $updated = $updated.Replace("Time:                        20:59:55", "Time:                        $newTime")
    $updated = $updated.Replace("Time:                        20:59:56", "Time:                        $newTime")

    $cell.Value = $updated
}
